$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: Task renamed, and Date Last Updated becomes a real date (04/03/2024),
#             reusing the date style already used elsewhere in column D (e.g. D4).
$ws.Range("A13").Value = "Task 13: Projects Page Frontend Polish (This Includes Projects, tasks, & members)"
$ws.Range("D4").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = 45385

# --- Row 18: Task 18 moves from Pending -> Complete, gains a Date Last Updated (text)
$ws.Range("B2").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Complete"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "03:27/2024"
$ws.Range("D14").Copy()
$ws.Range("D18").PasteSpecial(-4122)

# --- Row 19: Task 19 moves from Pending -> In Development, gains a Date Last Updated (text)
$ws.Range("G3").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "In Development"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "03/28/2024"
$ws.Range("D14").Copy()
$ws.Range("D19").PasteSpecial(-4122)

# --- Row 20: Date Last Updated style picks up the left/top-aligned date format
$ws.Range("D4").Copy()
$ws.Range("D20").PasteSpecial(-4122)

# --- Row 24: gains a Date Last Updated value (text)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "03/22/2024"
$ws.Range("D14").Copy()
$ws.Range("D24").PasteSpecial(-4122)

# --- Row 27: gains a Notes value
$ws.Range("C27").Value = "this requires auth, alternative is v-if ownerid=user.id, only for project owner"

# --- Row 28: Task renamed from the old "...Frontend Polish..." text to "Task 28: Bug Fixes"
$ws.Range("A28").Value = "Task 28: Bug Fixes"

# --- Row 29: Task 29 text trimmed (dropped "Storage (Local)")
$ws.Range("A29").Value = "Task 29: "

# --- Row 30: Task 30 text gains trailing spaces
$ws.Range("A30").Value = "Task 30:  "

# --- Selection / view state
$ws.Activate()
$ws.Range("C35").Select()
